$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price values so Excel does not
# reinterpret them (with "." as thousands separators) as numbers.
$textCells = @(
    "D5",
    "D7",
    "D9",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D19",
    "D22",
    "D24",
    "D26",
    "D27",
    "D30",
    "D35",
    "D40",
    "D43",
    "D44",
    "D45",
    "D46",
    "D48"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Updated coin price / volume figures
$ws.Range("D2").Value = "38.686.52"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "2.084.97"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "228.41"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("D7").Value = "59.89"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "2.393.91"
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").Value = "15.01"
$ws.Range("E13").Value = "  +4.00%  "
$ws.Range("D14").Value = "21.85"
$ws.Range("E14").Value = "  +2.00%  "
$ws.Range("D15").Value = "0.800"
$ws.Range("E15").Value = "  +4.74%  "
$ws.Range("D16").Value = "5.49"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "2.103.91"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "38.652.18"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "71.47"
$ws.Range("E19").Value = "  +2.91%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D22").Value = "226.91"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "2.38"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").Value = "171.00"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("E28").Value = "  +7.77%  "
$ws.Range("E29").Value = "  +13.34%  "
$ws.Range("D30").Value = "19.13"
$ws.Range("E30").Value = "  +1.82%  "
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("D35").Value = "0.0609"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "17.88"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("E41").Value = "  +5.66%  "
$ws.Range("D42").Value = "1.544.56"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D43").Value = "100.25"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "0.0925"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.81"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").Value = "7.71"
$ws.Range("E46").Value = "  +8.84%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").Value = "4.10"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "2.282.34"

# Restore default (General) style on the cells we force-formatted as text
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
